$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("D3").Value = "2016-01-07 14:16:41"
$ws.Range("G3").Value = "2016-01-07 14:17:25"

$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("D3").Value = "2016-01-07 14:16:53"
$ws2.Range("G3").Value = "2016-01-07 14:17:46"
